# Added Wins/Losses/Ties columns (AD:AF) populated with each team's season record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): copy the existing header formatting onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows (2-47): every row gets the same season record.
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 76
    $ws.Cells.Item($row, 31).Value = 86
    $ws.Cells.Item($row, 32).Value = 0
}
